$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.858.14"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.636.34"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'215.07"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.5022"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "'0.06374"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'19.70"
$ws.Range("D11").Value = "'0.07699"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "1.687.82"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "'4.260"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "1.860.89"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'0.5448"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "0.0₅7895"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "'64.27"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "25.856.47"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'202.96"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").Value = "'4.362"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "'9.896"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "'5.971"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").Value = "'1.920"
$ws.Range("E25").Value = "  +9.84%  "
$ws.Range("D26").Value = "'141.19"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'0.1135"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "'6.713"
$ws.Range("E29").Value = "  -3.78%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'0.04927"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "'3.188"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'1.539"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'2.626"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("D37").Value = "'0.8923"
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("D38").Value = "1.161.71"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.5600"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "'0.01560"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'5.709"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "'99.67"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.772.64"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").Value = "'0.4513"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'54.89"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  -0.53%  "
